$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlPasteFormats constant used below to copy number format / font / border
# from an existing, similarly-formatted cell onto the newly added ones.
$xlPasteFormats = -4122

# --- Row 4 (years header): add 2019 (P4) and 2020 (Q4), matching O4's style ---
$ws.Range("P4").Value = 2019
$ws.Range("O4").Copy() | Out-Null
$ws.Range("P4").PasteSpecial($xlPasteFormats)

$ws.Range("Q4").Value = 2020
$ws.Range("O4").Copy() | Out-Null
$ws.Range("Q4").PasteSpecial($xlPasteFormats)

# --- Row 5 (share of renewables, %): add 35.67 (P5) and a blank but styled Q5 ---
$ws.Range("P5").Value = 35.67
$ws.Range("E5").Copy() | Out-Null
$ws.Range("P5").PasteSpecial($xlPasteFormats)

$ws.Range("Q5").Value = $null
$ws.Range("E5").Copy() | Out-Null
$ws.Range("Q5").PasteSpecial($xlPasteFormats)

# --- Row 6 (hydropower production, mln kWh): add 13859.3 (P6) and 13979.1 (Q6) ---
$ws.Range("P6").Value = 13859.3
$ws.Range("O6").Copy() | Out-Null
$ws.Range("P6").PasteSpecial($xlPasteFormats)

$ws.Range("Q6").Value = 13979.1
$ws.Range("O6").Copy() | Out-Null
$ws.Range("Q6").PasteSpecial($xlPasteFormats)

# --- Restore the active selection to P9, as in the target workbook ---
$ws.Range("P9").Select() | Out-Null
